# Updated database table list
# Adds a new row (43) to the "MySQL" sheet describing the
# youshanding_hangzhou_claim_cmp analysis table, mirroring the existing
# rows' layout/formatting, and updates the view selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MySQL")

# --- New row 43 content --------------------------------------------------
$ws.Cells.Item(43, 1).Value = "youshanding_hangzhou_claim_cmp"
$ws.Cells.Item(43, 2).Value = "analysis"
$ws.Cells.Item(43, 3).Value = "分析"
$ws.Cells.Item(43, 4).Value = "优闪订数据与杭州本地化价格对比分析"
$ws.Cells.Item(43, 5).Value = "youshanding_hangzhou_cmp"
$ws.Cells.Item(43, 6).Value = "hangzhou_claim"
$ws.Cells.Item(43, 7).Value = "youshanding"
$ws.Cells.Item(43, 8).Value = 42949

# Match formatting used by the other data rows: wrapped text for B:G and
# the short-date display for the "creation time" column H.
$ws.Range("B43:G43").WrapText = $true
$ws.Cells.Item(43, 8).NumberFormat = "m/d/yy"

# Other wrapped rows in this sheet (e.g. row 42, the template this row was
# cloned from) render at this same taller row height.
$ws.Range("A43:H43").RowHeight = 40.5

# --- View state ------------------------------------------------------------
[void]$ws.Activate()
[void]$ws.Range("J42").Select()
